# Jess Midterm Subset Scores
# Adds a "Grade (/10)" column (computed from RMLSE), renames the Notes
# column to Model's notes-as-model-name, and appends six new rows for
# additional models (forward/backward stepwise variants), replacing the
# generic "KNN" notes with fuller descriptive model names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New numeric inputs (column A) for the added rows -----------------
$ws.Range("A4").Value = 0.29415
$ws.Range("A5").Value = 0.22727
$ws.Range("A6").Value = 0.22928
$ws.Range("A7").Value = 0.21879
$ws.Range("A8").Value = 0.23554
$ws.Range("A9").Value = 0.21734

# --- Column C (Model) header + new / existing rows ---------------------
# Re-assert the header text so the shared-string table is rebuilt in the
# same order the original author typed things in.
$ws.Range("C1").Value = "Model"

$ws.Range("C4").Value = "Forward stepwise AdjR2"
$ws.Range("C5").Value = "Forward stepwise BIC"
$ws.Range("C6").Value = "Forward stepwise Cp"
$ws.Range("C7").Value = "Backward stepwise AdjR2"
$ws.Range("C8").Value = "Backward stepwise BIC"
$ws.Range("C9").Value = "Backward stepwise Cp"

# --- Column B header + grade formula -----------------------------------
$ws.Range("B1").Value = "Grade (/10)"

# Reword the old generic "KNN" notes into the two rows that keep them.
$ws.Range("C2").Value = "KNN sqrt(n)/2 quantitative only"
$ws.Range("C3").Value = "KNN sqrt(n) quantitative only"

# B2 gets its own (non-shared) copy of the formula, B3:B9 share one group.
$ws.Range("B2").Formula = "=MAX(0.15, 0.12/A2*10)"
$ws.Range("B3:B9").Formula = "=MAX(0.15, 0.12/A3*10)"
$ws.Range("B2:B9").NumberFormat = "0.00"

# --- Columns D & E for the new rows -------------------------------------
$ws.Range("D4").Value = "Jess"
$ws.Range("D5").Value = "Jess"
$ws.Range("D6").Value = "Jess"
$ws.Range("D7").Value = "Jess"
$ws.Range("D8").Value = "Jess"
$ws.Range("D9").Value = "Jess"

$ws.Range("E4").Value = 43397
$ws.Range("E5").Value = 43397
$ws.Range("E6").Value = 43397
$ws.Range("E7").Value = 43397
$ws.Range("E8").Value = 43397
$ws.Range("E9").Value = 43397

# --- Cosmetic touch-ups to match the saved view -------------------------
$ws.Columns.Item(2).ColumnWidth = 9.7
$ws.Columns.Item(3).ColumnWidth = 26.3
$ws.Columns.Item(5).ColumnWidth = 9.7

$ws.Range("H6").Select() | Out-Null
